# Fixed bug in getTestOutput() function
# The isotope replicate d13C readings on the "Isotopes-D3" sheet had been
# entered into the wrong column (F vs G) for a handful of rows; this moves
# each value back to the correct column and refreshes the STDEVA formula
# that depends on them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Isotopes-D3")

# --- Row 16/17: value moved from G16 to F16, G17's STDEVA range corrected ---
$ws.Range("F16").Value = -33.868000000000002
$ws.Range("G16").ClearContents()
$ws.Range("G17").Formula = "=STDEVA(F34,F17,F18,F35,F36,F16)"

# --- Rows 31-34: values moved between F and G ---
$ws.Range("G31").Value = -30.753999999999998
$ws.Range("F31").ClearContents()

$ws.Range("G32").Value = -31.276
$ws.Range("F32").ClearContents()

$ws.Range("G33").Value = -31.673999999999999
$ws.Range("F33").ClearContents()

$ws.Range("F34").Value = -32.392000000000003
$ws.Range("G34").ClearContents()

# --- Rows 40-44: values moved between F and G; stale STDEVA on G43 dropped ---
$ws.Range("G40").Value = -31.684999999999999
$ws.Range("F40").ClearContents()

$ws.Range("G41").Value = -31.58
$ws.Range("F41").ClearContents()

$ws.Range("G42").Value = -30.893000000000001
$ws.Range("F42").ClearContents()

$ws.Range("G43").Value = -26.995999999999999
$ws.Range("F43").ClearContents()

$ws.Range("G44").Value = -27.643999999999998
$ws.Range("F44").ClearContents()

# --- Rows 49-51: values moved between F and G ---
$ws.Range("G49").Value = -28.56
$ws.Range("F49").ClearContents()

$ws.Range("G50").Value = -29.062999999999999
$ws.Range("F50").ClearContents()

$ws.Range("G51").Value = -29.067999999999998
$ws.Range("F51").ClearContents()

# --- Rows 75-77: values moved between F and G ---
$ws.Range("G75").Value = -30.667999999999999
$ws.Range("F75").ClearContents()

$ws.Range("G76").Value = -31.349999999999998
$ws.Range("F76").ClearContents()

$ws.Range("G77").Value = -30.332999999999998
$ws.Range("F77").ClearContents()

# --- Rows 101-103: values moved between F and G ---
$ws.Range("G101").Value = -29.654
$ws.Range("F101").ClearContents()

$ws.Range("G102").Value = -26.983000000000001
$ws.Range("F102").ClearContents()

$ws.Range("G103").Value = -29.425999999999998
$ws.Range("F103").ClearContents()

# Re-fit the touched rows so they don't pick up a stale autofit row height
$ws.Range("A16:A103").EntireRow.AutoFit() | Out-Null

# --- Clear the stale autofilter criterion (was filtering to "S-III-13")
# and unhide the rows it had hidden ---
$ws.ShowAllData()

# --- Selection / active-sheet bookkeeping: Isotopes-D3 becomes the active
# tab, with E11 selected (previously Conc-D3 was active) ---
$ws.Activate()
$ws.Range("E11").Select() | Out-Null
